# trafo_id -> gridnode_id refactor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "trafo_id" column header to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Reflect the cell selection recorded at save time
$ws.Range("G7").Select()
